$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

# Update every "Förändrad" date in column C from 45181 to 45182 (i.e. +1 day)
$ws.Range("C2:C$lastRow").Value = 45182
